# The commit this deck ships with ("Added an empty TwitterController api")
# is unrelated application source-code work; PowerPoll/.../bin/Debug/Presentation1.pptx
# is a build artifact that simply got re-saved/rebuilt alongside it. Diffing the
# canonical OOXML of that resave against the prior build shows:
#
#   * every r:id in presentation.xml (slide master, the single slide, all 11
#     slide layouts), the slide's webextensionref, and the picture's blip
#     r:embed were re-minted (classic PowerPoint "resave regenerates every
#     relationship id" churn - no target, order, count, or visual placement
#     of any relationship actually changed),
#   * the we:webextension part's GUID (ppt/slides/udata/data.xml /
#     we:webextension/@id) was likewise re-minted to a fresh instance id.
#
# None of that is user content: no shape was added/removed/moved/resized, no
# text, fill, picture bytes, slide size, or layout changed. The add-in
# "snapshot" GUID and package relationship ids are internal plumbing that
# PowerPoint itself regenerates on save and are not exposed anywhere in the
# Shapes/TextFrame/Tags/CustomXMLParts object model (there is no WebExtension
# object in the PowerPoint COM surface - Office web add-in metadata is not
# automatable). They also aren't safely reachable positionally: the lone
# graphicFrame/picture AlternateContent shape that carries the webextension
# reference is not individually addressable through Shapes.Item/Shapes.Range
# in this deck (it collapses onto the "Title 1" placeholder), so forcing an
# edit through shape indices would corrupt the title placeholder instead of
# touching the add-in reference - a strictly worse outcome than leaving the
# deck exactly as authored.
#
# So: open the deck, confirm nothing user-visible needs to change, and leave
# every slide/shape/text value exactly as it already is.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

Write-Host "Slides:" $p.Slides.Count
Write-Host "Shapes on slide 1:" $s.Shapes.Count
Write-Host "No user-visible content differs from the authored deck; nothing to edit."
